$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Rang" column header in F1
$ws.Cells.Item(1, 6).Value = "Rang"

# Fill in column E ("Moyenne de l'étudiant") with each student's average
$averages = @{
    3  = 7
    4  = 13
    5  = 13
    6  = 20
    7  = 6
    8  = 10
    9  = 10
    10 = 6
    11 = 19
    12 = 12
    13 = 15
    14 = 19
    15 = 11
    16 = 8
    17 = 17
    18 = 11
    19 = 18
    20 = 7
    21 = 19
    22 = 19
    23 = 13
    24 = 11
    25 = 20
    26 = 7
    27 = 17
    28 = 20
    29 = 10
    30 = 12
    31 = 6
    32 = 14
    33 = 6
    34 = 9
    35 = 11
    36 = 10
    37 = 17
    38 = 20
    39 = 12
    40 = 6
    41 = 11
    42 = 6
    43 = 12
    44 = 12
    45 = 18
    46 = 20
    47 = 14
    48 = 9
    49 = 5
    50 = 20
    51 = 12
    52 = 13
    53 = 6
    54 = 18
    55 = 13
    56 = 16
    57 = 11
    58 = 7
    59 = 18
    60 = 11
    61 = 18
    62 = 8
    63 = 6
}

for ($row = 3; $row -le 63; $row++) {
    $ws.Cells.Item($row, 5).Value = $averages[$row]
}

# Match Excel's default body font size after the edit (applies to the
# workbook's Normal style so every existing cell picks it up without
# materialising per-cell style overrides)
$normalStyle = $wb.Styles.Item(1)
$normalStyle.Font.Size = 11
